$d = $word.ActiveDocument

# Locate the target paragraphs by their distinctive text.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "MonoPamDataPhysicalExtensor.m*") {
        $startPara = $i
        break
    }
}
for ($i = $startPara; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "RobotPamCalculation.m*") {
        $endPara = $i + 2
        break
    }
}

$start = $d.Paragraphs.Item($startPara).Range.Start
$end = $d.Paragraphs.Item($endPara).Range.End
$r = $d.Range($start, $end)

$xml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7D8F314D" w14:textId="7BEEC17D" w:rsidR="00623BAD" w:rsidRDefault="00623BAD" w:rsidP="00623BAD">
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>MonoPamDataPhysicalExtensor.m</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> – This class was used for Lindie and Ben’s physical knee test experiment. It is the MonoPamData.m constructor, however it hard codes the length of the PAM, tendon, and </w:t>
  </w:r>
  <w:r>
    <w:t>air fittings, instead of coming up with a solution on its own. Specifically used for the knee extensor.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>MonoPamDataPhysicalFlexor.m</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> – This class was used for Lindie and Ben’s physical knee test experiment. It is the MonoPamData.m constructor, however it hard codes the length of the PAM, tendon, and air fittings, instead of coming up with a solution on its own. Specifically used for the knee flexor.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>RobotPamCalculation.m</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> -</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>This script contains all of the locations and other variables necessary for the PAM classes to calculate. The entire script can be ran to calculate all PAM values or specific PAM properties can be copied and pasted into optimization scripts. It also includes transformation matrices for the hip and knee joint.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:bookmarkStart w:id="0" w:name="_Hlk99733125"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>MonoPamDataP</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>inned</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Extensor.m</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t xml:space="preserve">– </w:t>
  </w:r>
  <w:r>
    <w:t>Same as</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> MonoPam</w:t>
  </w:r>
  <w:r>
    <w:t>PhysicalExtensor</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> constructor</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> but for the pinned knee joint.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>MonoPamDataPhysicalFlexor.m</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> – </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Same as </w:t>
  </w:r>
  <w:r>
    <w:t>MonoPamDataPinnedExtensor.m</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> but for a Flexor muscle.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>MonoPamData</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Explicit</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>.m</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>-  Will calculate everything given all the inputs:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>name, location, cross, diameter, t, rest, kmax, tendon, fit, pres</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>If only given the first 7 inputs its will assume tendon length = 0, pressure = 620 kPa, and fitting length = 0.0254 m.</w:t>
  </w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xml)
Write-Output "replaced paragraphs $startPara to $endPara"
